#
# Inserts a new "FY" / "Future Year" acronym row into the "Key to Variables"
# sheet, immediately above the existing "IT" / "Initial Time" row (which is
# the last row of the "plcy-schd" Top Level Folder group), shifting every
# row below it down by one.
#
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")
$ws.Activate()

# Insert a new row above row 209 (the thick-bottom-bordered "IT" row).
# Excel's Insert() copies formatting from the row above the insertion point,
# which already carries the correct (non-bordered) style for this new row.
$ws.Rows.Item(209).Insert()

# Fill in the new row's contents.
$ws.Cells.Item(209, 1).Value2 = "plcy-schd"
$ws.Cells.Item(209, 2).Value2 = "FY"
$ws.Cells.Item(209, 3).Value2 = "Future Year"
$ws.Cells.Item(209, 6).Value2 = "n/a"
$ws.Cells.Item(209, 7).Value2 = "If you change Initial Time, you need to update this file"

# Restore the view state: selection on the newly added note cell, and the
# frozen-pane viewport scrolled so row 191 is the first visible row below
# the frozen header.
$ws.Range("G209").Select()
$excel.ActiveWindow.ScrollRow = 191
